# Scheduled market-data refresh: updates plain (non-formula) price/profit
# values on the per-job "Leve profit" tables (columns H-N) across the
# ALC, ARM, BSM, CRP, GSM, LTW and WVR sheets. A few rows also gain or
# lose their LeveProfitNQ/LeveProfitHQ (M/N) cells entirely, matching
# rows whose computed profit became available/unavailable in this run.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$values = @{
    "H19" = 842.44446
    "J19" = 824.75
    "L19" = 824.75
    "N19" = -1174.75
    "H69" = 9237.700000000001
    "I69" = 4833
    "J69" = 10015
    "K69" = 14499
    "L69" = 30045
    "M69" = -13625
    "N69" = -31793
    "H72" = 9237.700000000001
    "I72" = 4833
    "J72" = 10015
    "K72" = 43497
    "L72" = 90135
    "M72" = -39129
    "N72" = -98871
    "H92" = 14286117
    "I92" = 17241828
    "J92" = 177.5
    "K92" = 17241828
    "L92" = 177.5
    "M92" = -17240580
    "N92" = -2673.5
    "H98" = 1381.9592
    "I98" = 1422.8864
    "K98" = 1422.8864
    "M98" = 75.11359999999991
    "H106" = 4170244.5
    "I106" = 5131639
    "K106" = 5131639
    "M106" = -5131008
    "H113" = 1943.64
    "I113" = 1645.8235
    "J113" = 2576.5
    "K113" = 1645.8235
    "L113" = 2576.5
    "M113" = 1608.1765
    "N113" = -9084.5
    "H116" = 15000.389
    "I116" = 17977.54
    "J116" = 7259.8
    "K116" = 17977.54
    "L116" = 7259.8
    "M116" = -14535.54
    "N116" = -14143.8
    "H122" = 1381.9592
    "I122" = 1422.8864
    "K122" = 4268.6592
    "M122" = -1818.6592
    "H131" = 5457.6665
    "I131" = 1054.75
    "K131" = 3164.25
    "M131" = 1875.75
    "H132" = 3812.1914
    "J132" = 2259.7144
    "L132" = 6779.1432
    "N132" = -11839.1432
    "H141" = 3507.7273
    "I141" = 2532
    "K141" = 7596
    "M141" = -2416
}
foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$values = @{
    "H32" = 1868.0605
    "I32" = 1910.5
    "K32" = 1910.5
    "M32" = -1623.5
    "H92" = 65665.164
    "J92" = 65665.164
    "L92" = 65665.164
    "N92" = -70657.164
    "H110" = 2252.037
    "I110" = 2205.92
    "K110" = 2205.92
    "M110" = -160.9200000000001
    "H119" = 89995
    "J119" = 89995
    "L119" = 89995
    "N119" = -99671
    "H122" = 4157.8604
    "I122" = 3231.7942
    "K122" = 9695.382599999999
    "M122" = -7245.382599999999
    "H132" = 3544.8262
    "I132" = 3037.1924
    "K132" = 9111.5772
    "M132" = -6581.5772
    "H138" = 0
    "I138" = 0
    "J138" = 0
    "K138" = 0
    "L138" = 0
    "H139" = 0
    "J139" = 0
    "L139" = 0
    "H140" = 0
    "J140" = 0
    "L140" = 0
    "H141" = 999999
    "J141" = 999999
    "L141" = 999999
    "N141" = -1010359
}
foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
$clearAddrs = @("M138", "N138", "N139", "N140")
foreach ($addr in $clearAddrs) {
    $ws.Range($addr).ClearContents()
}

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$values = @{
    "H86" = 3559.6
    "I86" = 2444.9546
    "J86" = 6624.875
    "K86" = 2444.9546
    "L86" = 6624.875
    "M86" = -1321.9546
    "N86" = -8870.875
    "H89" = 3559.6
    "I89" = 2444.9546
    "J89" = 6624.875
    "K89" = 12224.773
    "L89" = 33124.375
    "M89" = -6608.773000000001
    "N89" = -44356.375
    "H97" = 46345.285
    "I97" = 45142.668
    "J97" = 47247.25
    "K97" = 45142.668
    "L97" = 47247.25
    "M97" = -44151.668
    "N97" = -49229.25
    "H105" = 2173.6924
    "I105" = 1775.8
    "K105" = 1775.8
    "M105" = -28.79999999999995
    "H107" = 2001.4
    "I107" = 1251.8125
    "J107" = 4999.75
    "K107" = 1251.8125
    "L107" = 4999.75
    "M107" = 668.1875
    "N107" = -8839.75
}
foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$values = @{
    "H97" = 96385.336
    "J97" = 96385.336
    "L97" = 96385.336
    "N97" = -98367.336
}
foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$values = @{
    "H64" = 67000
    "J64" = 67000
    "L64" = 67000
    "N64" = -67496
    "H67" = 67000
    "J67" = 67000
    "L67" = 67000
    "N67" = -68716
    "H75" = 111849
    "J75" = 111849
    "L75" = 111849
    "N75" = -113597
    "H78" = 111849
    "J78" = 111849
    "L78" = 335547
    "N78" = -344283
    "H113" = 2190.25
    "I113" = 2092.9333
    "K113" = 2092.9333
    "M113" = 77.06669999999986
    "H122" = 3430.1052
    "I122" = 2974.2307
    "J122" = 4417.8335
    "K122" = 8922.6921
    "L122" = 13253.5005
    "M122" = -6472.6921
    "N122" = -18153.5005
    "H126" = 2892.6365
    "I126" = 2471.125
    "J126" = 4016.6667
    "K126" = 7413.375
    "L126" = 12050.0001
    "M126" = -4943.375
    "N126" = -16990.0001
    "H132" = 0
    "I132" = 0
    "J132" = 0
    "K132" = 0
    "L132" = 0
}
foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
$clearAddrs = @("M132", "N132")
foreach ($addr in $clearAddrs) {
    $ws.Range($addr).ClearContents()
}

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$values = @{
    "H40" = 4198.421
    "I40" = 4198.421
    "K40" = 4198.421
    "M40" = -4062.421
    "H61" = 1856
    "I61" = 1727.25
    "K61" = 1727.25
    "M61" = -1525.25
    "H113" = 1856
    "I113" = 1727.25
    "K113" = 1727.25
    "M113" = 442.75
    "H122" = 4435
    "I122" = 4150
    "K122" = 12450
    "M122" = -10000
    "H132" = 3716.1667
    "I132" = 3432.3333
    "K132" = 10296.9999
    "M132" = -7766.999899999999
    "H136" = 4522.222
    "I136" = 4125
    "J136" = 4840
    "K136" = 12375
    "L136" = 14520
    "M136" = -9825
    "N136" = -19620
}
foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$values = @{
    "H5" = 12666.667
    "J5" = 12666.667
    "L5" = 12666.667
    "N5" = -12890.667
    "H46" = 48642
    "J46" = 48642
    "L46" = 48642
    "N46" = -49104
    "H109" = 27233.334
    "J109" = 29800
    "L109" = 29800
    "N109" = -32574
    "H126" = 40276.383
    "I126" = 77084
    "J126" = 8727
    "K126" = 231252
    "L126" = 26181
    "M126" = -228782
    "N126" = -31121
    "H134" = 48642
    "J134" = 48642
    "L134" = 145926
    "N134" = -150996
    "H136" = 50963.285
    "I136" = 2620.4167
    "J136" = 115420.445
    "K136" = 7861.250100000001
    "L136" = 346261.335
    "M136" = -5311.250100000001
    "N136" = -351361.335
    "H138" = 112500
    "J138" = 112500
    "L138" = 112500
    "N138" = -122780
}
foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
